$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 changes
$ws.Range("L2").Value = 1.36
$ws.Range("M2").Value = 3

# Row 3 changes
$ws.Range("G3").Value = 3.1
$ws.Range("H3").Value = 2.75
$ws.Range("I3").Value = 2.63
$ws.Range("J3").Value = 1.1
$ws.Range("K3").Value = 7
$ws.Range("P3").Value = 1.5
$ws.Range("Q3").Value = 2.5
$ws.Range("R3").Value = 1.83
$ws.Range("S3").Value = 1.83
$ws.Range("U3").Value = 15
$ws.Range("Z3").Value = 7
$ws.Range("AE3").Value = 7.5
